# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Mon May  1 08:14:24 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.644.90"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.89%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.02"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.91%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4640"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.68%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3905"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.80%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.40"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07910"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9840"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.28"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -6.23%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.851.89"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.60%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.843"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.003"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06838"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.02%  "

$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.62"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001007"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.12"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.83%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.660.90"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.390"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.135"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.068.85"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.23"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.49"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.067"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.023"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.72"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9751"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09419"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.370"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.481"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.352"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06178"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02199"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.00%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.002"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5723"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.605"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.19"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1802"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.65%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.248"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5396"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.75"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07150"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.911"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "114.46"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.03%  "

